# Generate Report for Handback
# Adds a new handed-back file (d754c346-5338-4dfe-aee8-5240cb31b52f.md) as a
# fourth row to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$fileBase   = "d754c346-5338-4dfe-aee8-5240cb31b52f"
$mdName     = "$fileBase.md"
$mdPath     = "e2e\$fileBase.md"
$hashZh     = "6a1c5b9e4e3f0dfd24f4ae92cab7694f70a61c64"
$hashDe     = "6a1c5b9e4e3f0dfd24f4ae92cab7694f70a61c64"
$xlfZh      = "$fileBase.$hashZh.zh-cn.xlf"
$xlfDe      = "$fileBase.$hashDe.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$extMd        = ".md"
$dateFmt      = "yyyy-mm-dd HH:mm:ss"

$hyperColor = 15570276  # matches the workbook's custom hyperlink color FF6495ED

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOv = $rowOverview.Range

$rngOv.Item(1).Value = $mdName
$rngOv.Item(2).Value = $mdPath
$rngOv.Item(3).Value = $extMd
$rngOv.Item(5).Value = $statusInSync
$rngOv.Item(6).Value = $statusInSync
$rngOv.Item(7).Value = "2016-09-07 01:01:20"
$rngOv.Item(7).NumberFormat = $dateFmt

$ovLink = $wsOverview.Hyperlinks.Add($rngOv.Item(2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f08f72d32c09e57e156145056ddd1ac1ce6384c/$mdPath", "", "", $mdPath)
$rngOv.Item(2).Font.Color = $hyperColor
$rngOv.Item(2).Font.Underline = 2

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rZh = $rowZh.Range

$rZh.Item(1).Value = $mdName
$rZh.Item(2).Value = $extMd
$rZh.Item(3).Value = $statusInSync
$rZh.Item(4).Value = "e2e"
$rZh.Item(5).Value = "ht"
$rZh.Item(6).Value = "'True"
$rZh.Item(7).Value = $xlfZh
$rZh.Item(8).Value = "2016-09-07 01:01:14"
$rZh.Item(8).NumberFormat = $dateFmt
$rZh.Item(9).Value = $mdName
$rZh.Item(10).Value = $xlfZh
$rZh.Item(11).Value = "2016-09-07 01:01:32"
$rZh.Item(11).NumberFormat = $dateFmt
$rZh.Item(12).Value = "'"
$rZh.Item(13).Value = "'True"
$rZh.Item(14).Value = "'"
$rZh.Item(15).Value = "'False"
$rZh.Item(16).Value = "'"

$zhLinkA = $wsZh.Hyperlinks.Add($rZh.Item(1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f08f72d32c09e57e156145056ddd1ac1ce6384c/$mdName", "", "", $mdName)
$rZh.Item(1).Font.Color = $hyperColor
$rZh.Item(1).Font.Underline = 2

$zhLinkI = $wsZh.Hyperlinks.Add($rZh.Item(9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$hashZh/$mdName", "", "", $mdName)
$rZh.Item(9).Font.Color = $hyperColor
$rZh.Item(9).Font.Underline = 2

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rDe = $rowDe.Range

$rDe.Item(1).Value = $mdName
$rDe.Item(2).Value = $extMd
$rDe.Item(3).Value = $statusInSync
$rDe.Item(4).Value = "e2e"
$rDe.Item(5).Value = "ht"
$rDe.Item(6).Value = "'True"
$rDe.Item(7).Value = $xlfDe
$rDe.Item(8).Value = "2016-09-07 01:01:20"
$rDe.Item(8).NumberFormat = $dateFmt
$rDe.Item(9).Value = $mdName
$rDe.Item(10).Value = $xlfDe
$rDe.Item(11).Value = "2016-09-07 01:01:40"
$rDe.Item(11).NumberFormat = $dateFmt
$rDe.Item(12).Value = "'"
$rDe.Item(13).Value = "'True"
$rDe.Item(14).Value = "'"
$rDe.Item(15).Value = "'False"
$rDe.Item(16).Value = "'"

$deLinkA = $wsDe.Hyperlinks.Add($rDe.Item(1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f08f72d32c09e57e156145056ddd1ac1ce6384c/$mdName", "", "", $mdName)
$rDe.Item(1).Font.Color = $hyperColor
$rDe.Item(1).Font.Underline = 2

$deLinkI = $wsDe.Hyperlinks.Add($rDe.Item(9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$hashDe/$mdName", "", "", $mdName)
$rDe.Item(9).Font.Color = $hyperColor
$rDe.Item(9).Font.Underline = 2

Write-Host "Report row added to Overview, zh-cn and de-de sheets."
